$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 57 (ALC)
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 45000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 135000
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -135998

# Row 74 (ALC)
$ws.Range("H74").Value = 1992.7142
$ws.Range("I74").Value = 2069.8
$ws.Range("J74").Value = 1800
$ws.Range("K74").Value = 2069.8
$ws.Range("L74").Value = 1800
$ws.Range("M74").Value = -1133.8
$ws.Range("N74").Value = -3672

# Row 77 (ALC)
$ws.Range("H77").Value = 1992.7142
$ws.Range("I77").Value = 2069.8
$ws.Range("J77").Value = 1800
$ws.Range("K77").Value = 10349
$ws.Range("L77").Value = 9000
$ws.Range("M77").Value = -5669
$ws.Range("N77").Value = -18360

# Row 116 (ALC)
$ws.Range("H116").Value = 6125
$ws.Range("J116").Value = 6125
$ws.Range("L116").Value = 6125
$ws.Range("N116").Value = -13009

# Row 132 (ALC)
$ws.Range("H132").Value = 940.8889
$ws.Range("I132").Value = 937.6667
$ws.Range("J132").Value = 966.6667
$ws.Range("K132").Value = 2813.0001
$ws.Range("L132").Value = 2900.0001
$ws.Range("M132").Value = -283.0001000000002
$ws.Range("N132").Value = -7960.0001

# Row 135 (ALC)
$ws.Range("H135").Value = 594.1429000000001
$ws.Range("I135").Value = 432
$ws.Range("K135").Value = 3888
$ws.Range("M135").Value = -1353

# Row 137 (ALC)
$ws.Range("H137").Value = 2643.9048
$ws.Range("J137").Value = 2795.1765
$ws.Range("L137").Value = 8385.529500000001
$ws.Range("N137").Value = -13485.5295

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 4652762
$ws.Range("J2").Value = 1999
$ws.Range("L2").Value = 1999
$ws.Range("N2").Value = -2225

# Row 32 (ARM)
$ws.Range("H32").Value = 3483.6038
$ws.Range("I32").Value = 2798.641
$ws.Range("K32").Value = 2798.641
$ws.Range("M32").Value = -2511.641

# Row 61 (ARM)
$ws.Range("H61").Value = 3448.1
$ws.Range("I61").Value = 1746.8334
$ws.Range("K61").Value = 1746.8334
$ws.Range("M61").Value = -1534.8334

# Row 63 (ARM)
$ws.Range("H63").Value = 7499
$ws.Range("I63").Value = 7499
$ws.Range("K63").Value = 7499
$ws.Range("M63").Value = -6813

# Row 66 (ARM)
$ws.Range("H66").Value = 7499
$ws.Range("I66").Value = 7499
$ws.Range("K66").Value = 37495
$ws.Range("M66").Value = -34063

# Row 116 (ARM)
$ws.Range("H116").Value = 4652762
$ws.Range("J116").Value = 1999
$ws.Range("L116").Value = 1999
$ws.Range("N116").Value = -6587

# Row 132 (ARM)
$ws.Range("H132").Value = 3283.4333
$ws.Range("I132").Value = 3160.2
$ws.Range("J132").Value = 3899.6
$ws.Range("K132").Value = 9480.599999999999
$ws.Range("L132").Value = 11698.8
$ws.Range("M132").Value = -6950.599999999999
$ws.Range("N132").Value = -16758.8

# Row 136 (ARM)
$ws.Range("H136").Value = 3448.1
$ws.Range("I136").Value = 1746.8334
$ws.Range("K136").Value = 5240.5002
$ws.Range("M136").Value = -2690.5002

# Row 140 (ARM)
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 4652762
$ws.Range("J3").Value = 1999
$ws.Range("L3").Value = 1999
$ws.Range("N3").Value = -2227

# Row 134 (BSM)
$ws.Range("H134").Value = 2112.5
$ws.Range("I134").Value = 225
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 675
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = 1860
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
# Row 33 (CRP)
$ws.Range("H33").Value = 16169.4
$ws.Range("I33").Value = 12989.75
$ws.Range("J33").Value = 28888
$ws.Range("K33").Value = 12989.75
$ws.Range("L33").Value = 28888
$ws.Range("M33").Value = -12610.75
$ws.Range("N33").Value = -29646

# Row 58 (CRP)
$ws.Range("H58").Value = 4349722.5
$ws.Range("I58").Value = 7247914
$ws.Range("J58").Value = 2435
$ws.Range("K58").Value = 7247914
$ws.Range("L58").Value = 2435
$ws.Range("M58").Value = -7247711
$ws.Range("N58").Value = -2841

# Row 99 (CRP)
$ws.Range("H99").Value = 1668499.5
$ws.Range("I99").Value = 1668499.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1668499.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1667001.5
$ws.Range("N99").ClearContents()

# Row 126 (CRP)
$ws.Range("H126").Value = 1668499.5
$ws.Range("I126").Value = 1668499.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5005498.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5003028.5
$ws.Range("N126").ClearContents()

# Row 134 (CRP)
$ws.Range("H134").Value = 1339.1936
$ws.Range("I134").Value = 925.5769
$ws.Range("K134").Value = 2776.7307
$ws.Range("M134").Value = -241.7307000000001

# Row 136 (CRP)
$ws.Range("H136").Value = 4349722.5
$ws.Range("I136").Value = 7247914
$ws.Range("J136").Value = 2435
$ws.Range("K136").Value = 21743742
$ws.Range("L136").Value = 7305
$ws.Range("M136").Value = -21741192
$ws.Range("N136").Value = -12405

$ws = $wb.Worksheets.Item("CUL")
# Row 122 (CUL)
$ws.Range("H122").Value = 1082.6666
$ws.Range("I122").Value = 200.25
$ws.Range("K122").Value = 1802.25
$ws.Range("M122").Value = 647.75

$ws = $wb.Worksheets.Item("GSM")
# Row 126 (GSM)
$ws.Range("H126").Value = 1854591.2
$ws.Range("I126").Value = 2527904.2
$ws.Range("K126").Value = 7583712.600000001
$ws.Range("M126").Value = -7581242.600000001

# Row 132 (GSM)
$ws.Range("H132").Value = 1834516.4
$ws.Range("I132").Value = 2565936.5
$ws.Range("K132").Value = 7697809.5
$ws.Range("M132").Value = -7695279.5

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (LTW)
$ws.Range("H40").Value = 5514.357
$ws.Range("I40").Value = 2927.4546
$ws.Range("K40").Value = 2927.4546
$ws.Range("M40").Value = -2791.4546

# Row 122 (LTW)
$ws.Range("H122").Value = 10213.429
$ws.Range("I122").Value = 3750
$ws.Range("J122").Value = 12798.8
$ws.Range("K122").Value = 11250
$ws.Range("L122").Value = 38396.39999999999
$ws.Range("M122").Value = -8800
$ws.Range("N122").Value = -43296.39999999999

# Row 136 (LTW)
$ws.Range("H136").Value = 7249.4443
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 7892.143
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 23676.429
$ws.Range("M136").Value = -12450
$ws.Range("N136").Value = -28776.429

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 2719.8
$ws.Range("I81").Value = 2719.8
$ws.Range("K81").Value = 5439.6
$ws.Range("M81").Value = -4378.6

# Row 84 (WVR)
$ws.Range("H84").Value = 2719.8
$ws.Range("I84").Value = 2719.8
$ws.Range("K84").Value = 27198
$ws.Range("M84").Value = -21894

# Row 122 (WVR)
$ws.Range("H122").Value = 135518.5
$ws.Range("I122").Value = 135518.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 406555.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -404105.5
$ws.Range("N122").ClearContents()

# Row 123 (WVR)
$ws.Range("H123").Value = 37786.375
$ws.Range("J123").Value = 37786.375
$ws.Range("L123").Value = 37786.375
$ws.Range("N123").Value = -47586.375

# Row 125 (WVR)
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840

# Row 132 (WVR)
$ws.Range("H132").Value = 3528.2856
$ws.Range("I132").Value = 2733.3333
$ws.Range("J132").Value = 4124.5
$ws.Range("K132").Value = 8199.999899999999
$ws.Range("L132").Value = 12373.5
$ws.Range("M132").Value = -5669.999899999999
$ws.Range("N132").Value = -17433.5

# Row 136 (WVR)
$ws.Range("H136").Value = 29243266
$ws.Range("I136").Value = 50508468
$ws.Range("J136").Value = 3612.375
$ws.Range("K136").Value = 151525404
$ws.Range("L136").Value = 10837.125
$ws.Range("M136").Value = -151522854
$ws.Range("N136").Value = -15937.125

# Row 140 (WVR)
$ws.Range("H140").Value = 49485.4
$ws.Range("J140").Value = 49485.4
$ws.Range("L140").Value = 49485.4
$ws.Range("N140").Value = -59845.4

# Row 141 (WVR)
$ws.Range("H141").Value = 50802.145
$ws.Range("J141").Value = 50802.145
$ws.Range("L141").Value = 50802.145
$ws.Range("N141").Value = -61162.145
